$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'21"
$ws.Range("D2").Value = "'10"
$ws.Range("E2").Value = "'3"

$ws.Range("C3").Value = "'0"
$ws.Range("D3").Value = "'2"

$ws.Range("C4").Value = "'4"
$ws.Range("D4").Value = "'8"
$ws.Range("E4").Value = "'1"

$ws.Range("C5").Value = "'0"
$ws.Range("D5").Value = "'1"
$ws.Range("E5").Value = "'0"
$ws.Range("F5").Value = "'0"

$ws.Range("C6").Value = "'3"
$ws.Range("D6").Value = "'6"
$ws.Range("E6").Value = "'0"

$ws.Range("C7").Value = "'29"
$ws.Range("D7").Value = "'14"
$ws.Range("E7").Value = "'2"
$ws.Range("F7").Value = "'2"

$ws.Range("D8").Value = "'14"
$ws.Range("E8").Value = "'0"
